# Update "想去人数" (F column) values on the 展览 and 全部类型 sheets
# to reflect the latest scrape snapshot (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 14134
$ws1.Range("F3").Value  = 555
$ws1.Range("F5").Value  = 1218
$ws1.Range("F7").Value  = 13946
$ws1.Range("F8").Value  = 15047
$ws1.Range("F9").Value  = 6
$ws1.Range("F10").Value = 20
$ws1.Range("F12").Value = 179
$ws1.Range("F13").Value = 22
$ws1.Range("F20").Value = 22
$ws1.Range("F21").Value = 1165
$ws1.Range("F24").Value = 5843
$ws1.Range("F25").Value = 950
$ws1.Range("F26").Value = 1070
$ws1.Range("F27").Value = 5465
$ws1.Range("F28").Value = 59
$ws1.Range("F29").Value = 130
$ws1.Range("F30").Value = 70
$ws1.Range("F31").Value = 349

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 14134
$ws4.Range("F4").Value  = 555
$ws4.Range("F6").Value  = 1218
$ws4.Range("F8").Value  = 13946
$ws4.Range("F9").Value  = 15047
$ws4.Range("F10").Value = 6
$ws4.Range("F11").Value = 20
$ws4.Range("F13").Value = 179
$ws4.Range("F14").Value = 22
$ws4.Range("F21").Value = 22
$ws4.Range("F22").Value = 1165
$ws4.Range("F26").Value = 5843
$ws4.Range("F27").Value = 950
$ws4.Range("F28").Value = 1070
$ws4.Range("F29").Value = 5465
$ws4.Range("F30").Value = 59
$ws4.Range("F31").Value = 130
$ws4.Range("F32").Value = 70
$ws4.Range("F33").Value = 349

$wb.Save()
